# Updates the cryptos price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.264.04'
$ws.Range("E2").Value = '  +3.20%  '
$ws.Range("D3").Value = '3.063.99'
$ws.Range("E3").Value = '  +2.74%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''522.41'
$ws.Range("E5").Value = '  +3.87%  '
$ws.Range("D6").Value = '''142.49'
$ws.Range("E6").Value = '  +6.47%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +4.64%  '
$ws.Range("D9").Value = '''7.50'
$ws.Range("E9").Value = '  +2.58%  '
$ws.Range("E10").Value = '  +6.09%  '
$ws.Range("E11").Value = '  +5.57%  '
$ws.Range("D12").Value = '3.586.29'
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").Value = '''26.87'
$ws.Range("E14").Value = '  +6.73%  '
$ws.Range("D15").Value = '''0.0000171'
$ws.Range("E15").Value = '  +14.34%  '
$ws.Range("D16").Value = '58.234.55'
$ws.Range("E16").Value = '  +3.22%  '
$ws.Range("E17").Value = '  +9.45%  '
$ws.Range("D18").Value = '3.074.03'
$ws.Range("E18").Value = '  +3.16%  '
$ws.Range("E19").Value = '  +5.99%  '
$ws.Range("D20").Value = '''8.15'
$ws.Range("E20").Value = '  +4.91%  '
$ws.Range("D21").Value = '''338.59'
$ws.Range("E21").Value = '  +4.08%  '
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '''0.504'
$ws.Range("E24").Value = '  +6.96%  '
$ws.Range("D25").Value = '''65.47'
$ws.Range("D26").Value = '''0.169'
$ws.Range("E26").Value = '  +3.93%  '
$ws.Range("D27").Value = '0.0₃0964'
$ws.Range("E27").Value = '  +8.49%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = '''6.95'
$ws.Range("E29").Value = '  +7.32%  '
$ws.Range("D30").Value = '''7.56'
$ws.Range("E30").Value = '  +11.65%  '
$ws.Range("E31").Value = '  +5.72%  '
$ws.Range("E32").Value = '  +4.35%  '
$ws.Range("E33").Value = '  +4.60%  '
$ws.Range("E34").Value = '  +7.99%  '
$ws.Range("D35").Value = '''157.32'
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("E36").Value = '  +7.95%  '
$ws.Range("E37").Value = '  +2.51%  '
$ws.Range("E39").Value = '  +3.45%  '
$ws.Range("D40").Value = '3.101.14'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").Value = '''37.78'
$ws.Range("E41").Value = '  +4.58%  '
$ws.Range("D42").Value = '''3.92'
$ws.Range("E42").Value = '  +9.98%  '
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '2.336.48'
$ws.Range("E45").Value = '  +4.33%  '
$ws.Range("E46").Value = '  +4.71%  '
$ws.Range("E47").Value = '  +3.04%  '
$ws.Range("E48").Value = '  +5.81%  '
$ws.Range("E49").Value = '  +3.19%  '
$ws.Range("E50").Value = '  +5.10%  '
$ws.Range("D51").Value = '''1.88'
$ws.Range("E51").Value = '  -2.00%  '
